$wb = $excel.ActiveWorkbook

# --- Sheet "Confirmed" (sheet1): add row 29 for 2020-04-04 ---
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsConfirmed.Activate()
$wsConfirmed.Range("A28:C28").Copy() | Out-Null
$wsConfirmed.Range("A29:C29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsConfirmed.Cells.Item(29, 1).Value = 43925
$wsConfirmed.Cells.Item(29, 2).Value = 70
$wsConfirmed.Cells.Item(29, 3).Value = 9
$wsConfirmed.Range("D29").Select() | Out-Null

# --- Sheet "Recoverd" (sheet2): add row 29 for 2020-04-04 ---
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Activate()
$wsRecoverd.Range("A28:C28").Copy() | Out-Null
$wsRecoverd.Range("A29:C29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsRecoverd.Cells.Item(29, 1).Value = 43925
$wsRecoverd.Cells.Item(29, 2).Value = 30
$wsRecoverd.Cells.Item(29, 3).Value = 5
$wsRecoverd.Range("D29").Select() | Out-Null

# --- Sheet "Death" (sheet3): add row 29 for 2020-04-04 ---
$wsDeath = $wb.Worksheets.Item("Death")
$wsDeath.Activate()
$wsDeath.Range("A28:C28").Copy() | Out-Null
$wsDeath.Range("A29:C29").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsDeath.Cells.Item(29, 1).Value = 43925
$wsDeath.Cells.Item(29, 2).Value = 8
$wsDeath.Cells.Item(29, 3).Value = 2
$wsDeath.Range("B29").Select() | Out-Null

# --- Final view state: "Recoverd" tab is the active/selected tab ---
$wsRecoverd.Activate()
$wsRecoverd.Range("D29").Select() | Out-Null
